$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update store codes: row 19 LA -> N5, rows 20/21 NO_MAP -> XX / YY
$ws.Range("A19").Value = "N5"
$ws.Range("A20").Value = "ΧΧ"
$ws.Range("A21").Value = "ΥΥ"

# Rows 20 and 21: remove yellow highlight from columns B and C (keep it on A)
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Style = "Normal"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Style = "Normal"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 5.7109375
$ws.Columns.Item(2).ColumnWidth = 125

# Selection
$ws.Range("F9").Select() | Out-Null
